# Applies the Jan 12 2023 17:25 UTC symbol-list refresh described in the commit:
#  - Coin table rows 2..51 all move their "Hora" (hour) column from 16 -> 17
#  - Most rows also get refreshed Price / Volume(1h) figures
#  - Rows 7..19 shift: a new row (GateToken) is inserted at the top of that block,
#    pushing FTXToken..LEO down by one row each (GateToken's old slot at row 19
#    now holds LEO's refreshed data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text cells (Coin name / Link columns): safe to assign directly ---
$plainEdits = @(
    @{ Cell = 'B7'; Value = 'GateToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'B8'; Value = 'FTXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'B9'; Value = 'MXToken' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'B11'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'B13'; Value = 'BitrueCoin' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'B14'; Value = 'BitMartToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'B15'; Value = 'BitForexToken' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'B16'; Value = 'CoinExToken' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' }
    @{ Cell = 'B17'; Value = 'One' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one' }
    @{ Cell = 'B18'; Value = 'TigerCash' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'B19'; Value = 'LEO' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
)
foreach ($edit in $plainEdits) {
    $ws.Range($edit.Cell).Value = $edit.Value
}

# --- Numeric-looking text cells (Price / Volume(1h) / Hora columns) ---
# These must stay text (matching the original inline-string cells), so force
# NumberFormat to "@" (Text) before assigning, otherwise Excel's COM layer
# auto-converts strings like "282.53" or "17" into numeric values.
$textEdits = @(
    @{ Cell = 'D2'; Value = '282.53' }
    @{ Cell = 'E2'; Value = '2.00%' }
    @{ Cell = 'G2'; Value = '17' }
    @{ Cell = 'D3'; Value = '28.31' }
    @{ Cell = 'E3'; Value = '3.81%' }
    @{ Cell = 'G3'; Value = '17' }
    @{ Cell = 'D4'; Value = '5.046' }
    @{ Cell = 'E4'; Value = '4.23%' }
    @{ Cell = 'G4'; Value = '17' }
    @{ Cell = 'D5'; Value = '0.06508' }
    @{ Cell = 'E5'; Value = '2.93%' }
    @{ Cell = 'G5'; Value = '17' }
    @{ Cell = 'D6'; Value = '7.268' }
    @{ Cell = 'E6'; Value = '3.46%' }
    @{ Cell = 'G6'; Value = '17' }
    @{ Cell = 'D7'; Value = '3.367' }
    @{ Cell = 'E7'; Value = '2.07%' }
    @{ Cell = 'G7'; Value = '17' }
    @{ Cell = 'D8'; Value = '1.361' }
    @{ Cell = 'E8'; Value = '-0.15%' }
    @{ Cell = 'G8'; Value = '17' }
    @{ Cell = 'D9'; Value = '0.9291' }
    @{ Cell = 'E9'; Value = '4.58%' }
    @{ Cell = 'G9'; Value = '17' }
    @{ Cell = 'D10'; Value = '0.1558' }
    @{ Cell = 'E10'; Value = '3.22%' }
    @{ Cell = 'G10'; Value = '17' }
    @{ Cell = 'D11'; Value = '0.06190' }
    @{ Cell = 'E11'; Value = '13.35%' }
    @{ Cell = 'G11'; Value = '17' }
    @{ Cell = 'D12'; Value = '0.07620' }
    @{ Cell = 'E12'; Value = '2.01%' }
    @{ Cell = 'G12'; Value = '17' }
    @{ Cell = 'D13'; Value = '0.02904' }
    @{ Cell = 'E13'; Value = '0.43%' }
    @{ Cell = 'G13'; Value = '17' }
    @{ Cell = 'D14'; Value = '0.08967' }
    @{ Cell = 'E14'; Value = '0.19%' }
    @{ Cell = 'G14'; Value = '17' }
    @{ Cell = 'D15'; Value = '0.001593' }
    @{ Cell = 'E15'; Value = '1.38%' }
    @{ Cell = 'G15'; Value = '17' }
    @{ Cell = 'D16'; Value = '0.04457' }
    @{ Cell = 'E16'; Value = '1.87%' }
    @{ Cell = 'G16'; Value = '17' }
    @{ Cell = 'D17'; Value = '0.0006374' }
    @{ Cell = 'E17'; Value = '0.69%' }
    @{ Cell = 'G17'; Value = '17' }
    @{ Cell = 'D18'; Value = '0.006040' }
    @{ Cell = 'E18'; Value = '0.14%' }
    @{ Cell = 'G18'; Value = '17' }
    @{ Cell = 'D19'; Value = '3.451' }
    @{ Cell = 'E19'; Value = '-0.70%' }
    @{ Cell = 'G19'; Value = '17' }
    @{ Cell = 'D20'; Value = '2.235' }
    @{ Cell = 'E20'; Value = '0.09%' }
    @{ Cell = 'G20'; Value = '17' }
    @{ Cell = 'D21'; Value = '0.3203' }
    @{ Cell = 'E21'; Value = '1.03%' }
    @{ Cell = 'G21'; Value = '17' }
    @{ Cell = 'E22'; Value = '-4.62%' }
    @{ Cell = 'G22'; Value = '17' }
    @{ Cell = 'D23'; Value = '4.098' }
    @{ Cell = 'E23'; Value = '4.54%' }
    @{ Cell = 'G23'; Value = '17' }
    @{ Cell = 'D24'; Value = '0.1525' }
    @{ Cell = 'E24'; Value = '1.19%' }
    @{ Cell = 'G24'; Value = '17' }
    @{ Cell = 'D25'; Value = '0.001180' }
    @{ Cell = 'E25'; Value = '0.38%' }
    @{ Cell = 'G25'; Value = '17' }
    @{ Cell = 'D26'; Value = '0.004393' }
    @{ Cell = 'E26'; Value = '3.50%' }
    @{ Cell = 'G26'; Value = '17' }
    @{ Cell = 'D27'; Value = '0.0001247' }
    @{ Cell = 'E27'; Value = '5.76%' }
    @{ Cell = 'G27'; Value = '17' }
    @{ Cell = 'D28'; Value = '0.0001615' }
    @{ Cell = 'E28'; Value = '-2.07%' }
    @{ Cell = 'G28'; Value = '17' }
    @{ Cell = 'G29'; Value = '17' }
    @{ Cell = 'G30'; Value = '17' }
    @{ Cell = 'G31'; Value = '17' }
    @{ Cell = 'G32'; Value = '17' }
    @{ Cell = 'G33'; Value = '17' }
    @{ Cell = 'G34'; Value = '17' }
    @{ Cell = 'G35'; Value = '17' }
    @{ Cell = 'G36'; Value = '17' }
    @{ Cell = 'G37'; Value = '17' }
    @{ Cell = 'G38'; Value = '17' }
    @{ Cell = 'G39'; Value = '17' }
    @{ Cell = 'D40'; Value = '0.04158' }
    @{ Cell = 'E40'; Value = '4.24%' }
    @{ Cell = 'G40'; Value = '17' }
    @{ Cell = 'D41'; Value = '0.006601' }
    @{ Cell = 'E41'; Value = '-1.09%' }
    @{ Cell = 'G41'; Value = '17' }
    @{ Cell = 'D42'; Value = '0.1221' }
    @{ Cell = 'E42'; Value = '-12.47%' }
    @{ Cell = 'G42'; Value = '17' }
    @{ Cell = 'D43'; Value = '0.002016' }
    @{ Cell = 'E43'; Value = '-6.63%' }
    @{ Cell = 'G43'; Value = '17' }
    @{ Cell = 'D44'; Value = '0.01205' }
    @{ Cell = 'E44'; Value = '3.40%' }
    @{ Cell = 'G44'; Value = '17' }
    @{ Cell = 'D45'; Value = '0.00005529' }
    @{ Cell = 'E45'; Value = '-0.31%' }
    @{ Cell = 'G45'; Value = '17' }
    @{ Cell = 'G46'; Value = '17' }
    @{ Cell = 'D47'; Value = '0.01297' }
    @{ Cell = 'E47'; Value = '-29.86%' }
    @{ Cell = 'G47'; Value = '17' }
    @{ Cell = 'G48'; Value = '17' }
    @{ Cell = 'G49'; Value = '17' }
    @{ Cell = 'G50'; Value = '17' }
    @{ Cell = 'G51'; Value = '17' }
)
foreach ($edit in $textEdits) {
    $rng = $ws.Range($edit.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $edit.Value
}
